# Update "paises" (countries) workbook:
#  - Refresh the "Datos actualizados..." timestamp on A1
#  - Re-rank a few countries whose totals were updated (Irlanda overtakes
#    India/Ecuador, Camerun overtakes Kazajistan, Santa Lucia overtakes Sudan)
#    by swapping the country names/labels in column A so the table stays
#    sorted by "Casos totales" (column B) descending
#  - Update the numeric stats (columns B:H) for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 21:22"

# --- Re-sorted country labels (column A) -------------------------------
# Irlanda's updated totals now exceed India's and Ecuador's
$ws.Range("A24").Value = "Irlanda"
$ws.Range("A25").Value = "India"
$ws.Range("A26").Value = "Ecuador"

# Camerun's updated totals now exceed Kazajistan's
$ws.Range("A75").Value = "Camerun"
$ws.Range("A76").Value = "Kazajistan"

# Santa Lucia's updated totals now exceed Sudan's
$ws.Range("A178").Value = "Santa Lucia"
$ws.Range("A179").Value = "Sudan"

# --- Row 8 (Alemania) ---------------------------------------------------
$ws.Range("B8").Value = 120157
$ws.Range("C8").Value = 1922
$ws.Range("E8").Value = 65062
$ws.Range("G8").Value = 81
$ws.Range("H8").Value = 2688

# --- Row 19 (Corea del Sur) ----------------------------------------------
$ws.Range("B19").Value = 13555
$ws.Range("C19").Value = 311
$ws.Range("E19").Value = 7172

# --- Row 24 (now Irlanda, new updated numbers) ---------------------------
$ws.Range("B24").Value = 8089
$ws.Range("C24").Value = 1515
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 7777
$ws.Range("F24").Value = 194
$ws.Range("G24").Value = 24
$ws.Range("H24").Value = 287

# --- Row 25 (now India, keeps India's former numbers) ---------------------
$ws.Range("B25").Value = 7598
$ws.Range("C25").Value = 873
$ws.Range("D25").Value = 774
$ws.Range("E25").Value = 6578
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 246

# --- Row 26 (now Ecuador, keeps Ecuador's former numbers) -----------------
$ws.Range("B26").Value = 7161
$ws.Range("C26").Value = 2196
$ws.Range("D26").Value = 368
$ws.Range("E26").Value = 6496
$ws.Range("F26").Value = 171
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 297

# --- Row 53 (Bosnia y Herzegovina) ----------------------------------------
$ws.Range("E53").Value = 1650
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 92

# --- Row 73 (Armenia) -----------------------------------------------------
$ws.Range("B73").Value = 925
$ws.Range("C73").Value = 38
$ws.Range("D73").Value = 539
$ws.Range("E73").Value = 380

# --- Row 75 (now Camerun, new updated numbers) ----------------------------
$ws.Range("B75").Value = 820
$ws.Range("C75").Value = 17
$ws.Range("D75").Value = 98
$ws.Range("E75").Value = 710
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 12

# --- Row 76 (now Kazajistan, keeps Kazajistan's former numbers) -----------
$ws.Range("B76").Value = 812
$ws.Range("C76").Value = 31
$ws.Range("D76").Value = 64
$ws.Range("E76").Value = 738
$ws.Range("F76").Value = 21
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 10

# --- Row 88 (Afganistan) ---------------------------------------------------
$ws.Range("B88").Value = 558
$ws.Range("C88").Value = 19
$ws.Range("D88").Value = 42
$ws.Range("E88").Value = 513

# --- Row 161 (Eritrea) ------------------------------------------------------
$ws.Range("B161").Value = 31
$ws.Range("C161").Value = 1
$ws.Range("E161").Value = 29

# --- Row 178 (now Santa Lucia, updated numbers) -----------------------------
$ws.Range("C178").Value = 1
$ws.Range("D178").Value = 1
$ws.Range("E178").Value = 14
$ws.Range("H178").Value = 0

# --- Row 179 (now Sudan, updated numbers) ------------------------------------
$ws.Range("B179").Value = 15
$ws.Range("D179").Value = 2
$ws.Range("E179").Value = 11
$ws.Range("H179").Value = 0
